$wb = $excel.ActiveWorkbook

# Rename the metadata sheet 'General' to 'Table'
$ws = $wb.Worksheets.Item("General")
$ws.Name = "Table"

# Make this sheet the active / selected sheet (moves tabSelected + activeTab)
$ws.Activate()
